# Update forecasted energy values in column B (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78150.1224852913
$ws.Range("B3").Value = 73931.51120751112
$ws.Range("B4").Value = 70665.639883871
$ws.Range("B5").Value = 68221.08144896063
$ws.Range("B6").Value = 67442.62655372635
$ws.Range("B7").Value = 67697.40905076484
$ws.Range("B8").Value = 66805.99719848631
$ws.Range("B9").Value = 73518.8360027173
$ws.Range("B10").Value = 90665.86311579018
$ws.Range("B11").Value = 99986.44599549515
$ws.Range("B12").Value = 104010.8232173354
$ws.Range("B13").Value = 106016.5690550003
$ws.Range("B14").Value = 107320.0925482591
$ws.Range("B15").Value = 111547.1563226645
$ws.Range("B16").Value = 112002.2930044165
$ws.Range("B17").Value = 109747.5369375676
$ws.Range("B18").Value = 103869.0857919438
$ws.Range("B19").Value = 96041.01666867564
$ws.Range("B20").Value = 95816.15130808466
$ws.Range("B21").Value = 92433.58920697204
$ws.Range("B22").Value = 90165.2556892329
$ws.Range("B23").Value = 89664.82647591364
$ws.Range("B24").Value = 86158.80010198322
$ws.Range("B25").Value = 81709.70824788923
